# Auto-generated edit script for IFRS data correction
# Commit message: error solve ifrs list
#
# Row 2 had stale/misaligned financial data; it is cleared down to just
# AG/AH/AJ (which get refreshed values). Rows 3-9 get corrected values
# across nearly every numeric column, and Y3/Z3 are cleared.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$setMap = [ordered]@{
    "AG2" = 1000
    "AH2" = 3.58
    "AJ2" = 35862119
    "D3" = 25267
    "E3" = 2223
    "F3" = 2223
    "G3" = 2151
    "H3" = 1650
    "I3" = 1509
    "J3" = 141
    "K3" = 35939
    "L3" = 16773
    "M3" = 19166
    "N3" = 14732
    "O3" = 4434
    "P3" = 2193
    "Q3" = 3537
    "R3" = -4424
    "S3" = -560
    "T3" = 817
    "U3" = 2721
    "V3" = 9751
    "W3" = 8.800000000000001
    "X3" = 6.53
    "AA3" = 87.52
    "AB3" = 640.3
    "AC3" = 4208
    "AD3" = 6.56
    "AE3" = 41386
    "AF3" = 0.67
    "AG3" = 1000
    "AH3" = 3.62
    "AI3" = 23.59
    "AJ3" = 35862119
    "D4" = 25311
    "E4" = 1435
    "F4" = 1435
    "G4" = 1328
    "H4" = 1047
    "I4" = 887
    "J4" = 160
    "K4" = 35371
    "L4" = 15428
    "M4" = 19942
    "N4" = 17378
    "O4" = 2565
    "P4" = 2193
    "Q4" = 2552
    "R4" = -1714
    "S4" = -1791
    "T4" = 1608
    "U4" = 945
    "V4" = 9993
    "W4" = 5.67
    "X4" = 4.13
    "Y4" = 5.53
    "Z4" = 2.94
    "AA4" = 77.36
    "AB4" = 667.24
    "AC4" = 2473
    "AD4" = 10.33
    "AE4" = 48820
    "AF4" = 0.52
    "AG4" = 800
    "AH4" = 3.13
    "AI4" = 32.1
    "AJ4" = 35862119
    "D5" = 30553
    "E5" = 1885
    "F5" = 1885
    "G5" = 1733
    "H5" = 1372
    "I5" = 1267
    "J5" = 104
    "K5" = 35910
    "L5" = 14878
    "M5" = 21032
    "N5" = 18363
    "O5" = 2669
    "P5" = 2193
    "Q5" = 2219
    "R5" = -1255
    "S5" = -1043
    "T5" = 965
    "U5" = 1254
    "V5" = 9244
    "W5" = 6.17
    "X5" = 4.49
    "Y5" = 7.09
    "Z5" = 3.85
    "AA5" = 70.73999999999999
    "AB5" = 711.65
    "AC5" = 3533
    "AD5" = 8.220000000000001
    "AE5" = 51588
    "AF5" = 0.5600000000000001
    "AG5" = 900
    "AH5" = 3.1
    "AI5" = 25.28
    "AJ5" = 35862119
    "D6" = 32781
    "E6" = 559
    "F6" = 559
    "G6" = 295
    "H6" = 248
    "I6" = 210
    "K6" = 37189
    "L6" = 17239
    "M6" = 19950
    "N6" = 18596
    "P6" = 2193
    "Q6" = 880
    "R6" = -918
    "S6" = 631
    "T6" = 821
    "U6" = 59
    "V6" = 11498
    "W6" = 1.7
    "X6" = 0.76
    "Y6" = 1.14
    "Z6" = 0.68
    "AA6" = 86.41
    "AB6" = 708.58
    "AC6" = 585
    "AD6" = 30.24
    "AE6" = 52244
    "AF6" = 0.34
    "AG6" = 850
    "AH6" = 4.8
    "AI6" = 144.13
    "AJ6" = 35862119
    "D7" = 30021
    "E7" = 478
    "G7" = 361
    "H7" = 252
    "I7" = 238
    "K7" = 36515
    "L7" = 16644
    "M7" = 19871
    "N7" = 18513
    "P7" = 2191
    "Q7" = 2641
    "R7" = -1040
    "S7" = -1327
    "T7" = 1015
    "U7" = 1570
    "W7" = 1.59
    "X7" = 0.84
    "Y7" = 1.28
    "Z7" = 0.68
    "AA7" = 83.76000000000001
    "AC7" = 663
    "AD7" = 19.45
    "AE7" = 52009
    "AF7" = 0.25
    "AG7" = 730
    "AH7" = 5.66
    "AI7" = 110.04
    "D8" = 30024
    "E8" = 728
    "G8" = 579
    "H8" = 433
    "I8" = 389
    "K8" = 36451
    "L8" = 16422
    "M8" = 20029
    "N8" = 18640
    "P8" = 2191
    "Q8" = 2040
    "R8" = -838
    "S8" = -852
    "T8" = 799
    "U8" = 1255
    "W8" = 2.42
    "X8" = 1.44
    "Y8" = 2.1
    "Z8" = 1.19
    "AA8" = 81.98999999999999
    "AC8" = 1086
    "AD8" = 11.88
    "AE8" = 52367
    "AF8" = 0.25
    "AG8" = 750
    "AH8" = 5.81
    "AI8" = 69.08
    "D9" = 30794
    "E9" = 820
    "G9" = 674
    "H9" = 511
    "I9" = 463
    "K9" = 37053
    "L9" = 16793
    "M9" = 20260
    "N9" = 18836
    "P9" = 2191
    "Q9" = 2016
    "R9" = -989
    "S9" = -448
    "T9" = 894
    "U9" = 1213
    "W9" = 2.66
    "X9" = 1.66
    "Y9" = 2.47
    "Z9" = 1.39
    "AA9" = 82.88
    "AC9" = 1291
    "AD9" = 10
    "AE9" = 52918
    "AF9" = 0.24
    "AG9" = 760
    "AH9" = 5.89
    "AI9" = 58.89
}

$clearList = @(
    "D2"
    "E2"
    "F2"
    "G2"
    "H2"
    "I2"
    "J2"
    "K2"
    "L2"
    "M2"
    "N2"
    "O2"
    "P2"
    "Q2"
    "R2"
    "S2"
    "T2"
    "U2"
    "V2"
    "W2"
    "X2"
    "Y2"
    "Z2"
    "AA2"
    "AB2"
    "AC2"
    "AD2"
    "AE2"
    "AF2"
    "AI2"
    "Y3"
    "Z3"
)

foreach ($ref in $setMap.Keys) {
    $ws.Range($ref).Value = $setMap[$ref]
}

foreach ($ref in $clearList) {
    $ws.Range($ref).ClearContents()
}
